{"js": "// Highlight quantitative impact metrics (percentages, dollar amounts, and\n// large numbers) in specific resume bullet paragraphs using a hybrid\n// bold + color (#2C3E50) run format, matching the \"quantitative metrics\n// highlighting\" feature described in the commit message.\n\n// Map of the exact paragraph text -> ordered list of substrings within that\n// paragraph which must become bold + colored (2C3E50). Using the full\n// paragraph text as the lookup key keeps the edit scoped precisely to the\n// intended bullets (other paragraphs that happen to contain similar numbers,\n// e.g. the professional summary, are left untouched).\nconst HIGHLIGHTS = [\n  {\n    text: \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    tokens: [\"23%\", \"64%\"]\n  },\n  {\n    text: \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \\u00B14.2% to \\u00B12.1%\",\n    tokens: [\"87%\", \"71%\", \"\\u00B14.2%\", \"\\u00B12.1%\"]\n  },\n  {\n    text: \"\\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    tokens: [\"1,200\"]\n  },\n  {\n    text: \"\\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    tokens: [\"$400M\", \"$1B\"]\n  },\n  {\n    text: \"\\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    tokens: [\"73.5%\", \"$4.7M\"]\n  },\n  {\n    text: \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    tokens: [\"87%\", \"71%\"]\n  }\n];\n\nconst HIGHLIGHT_COLOR = \"2C3E50\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Track which HIGHLIGHTS entries have already been consumed so that the two\n// paragraphs sharing identical prefixes (\"Achieved 87% ... 71%\" appears both\n// with and without the trailing clause) are each matched to the correct\n// (and only the correct) paragraph instance, in document order.\nconst used = new Array(HIGHLIGHTS.length).fill(false);\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const paraText = para.text;\n\n  for (let h = 0; h < HIGHLIGHTS.length; h++) {\n    if (used[h]) continue;\n    if (paraText !== HIGHLIGHTS[h].text) continue;\n    used[h] = true;\n\n    for (const token of HIGHLIGHTS[h].tokens) {\n      const hits = para.search(token, { matchCase: true });\n      hits.load(\"items\");\n      await context.sync();\n\n      for (let k = 0; k < hits.items.length; k++) {\n        const hit = hits.items[k];\n        hit.font.bold = true;\n        hit.font.color = HIGHLIGHT_COLOR;\n      }\n      await context.sync();\n    }\n    break;\n  }\n}\n", "ps1": "# Highlight quantitative impact metrics (percentages, dollar amounts, and\n# large numbers) in specific resume bullet paragraphs using a hybrid\n# bold + color (#2C3E50) run format, matching the \"quantitative metrics\n# highlighting\" feature described in the commit message.\n\n$d = $word.ActiveDocument\n\n# Word's wdColor values are packed as 0x00BBGGRR (blue/green/red byte\n# order), not the usual 0xRRGGBB \u2014 convert the hex color used across the\n# other output formats (PDF/DOCX bold+color) into that integer form.\nfunction ConvertTo-WdColor([string]$hex) {\n    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)\n    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)\n    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)\n    return ($b * 65536) + ($g * 256) + $r\n}\n\n$highlightColor = ConvertTo-WdColor \"2C3E50\"\n\n# Map of the exact paragraph text -> ordered list of substrings within that\n# paragraph which must become bold + colored. Using the full paragraph text\n# as the lookup key keeps the edit scoped precisely to the intended bullets\n# (other paragraphs that happen to contain similar numbers, e.g. the\n# professional summary, are left untouched).\n$highlights = @(\n    @{\n        Text   = \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n        Tokens = @(\"23%\", \"64%\")\n    },\n    @{\n        Text   = \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\"\n        Tokens = @(\"87%\", \"71%\", \"\u00b14.2%\", \"\u00b12.1%\")\n    },\n    @{\n        Text   = \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\"\n        Tokens = @(\"1,200\")\n    },\n    @{\n        Text   = \"\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\"\n        Tokens = @(\"`$400M\", \"`$1B\")\n    },\n    @{\n        Text   = \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\"\n        Tokens = @(\"73.5%\", \"`$4.7M\")\n    },\n    @{\n        Text   = \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\"\n        Tokens = @(\"87%\", \"71%\")\n    }\n)\n\n# Snapshot matching paragraph indices up front (1-based, COM-style) before\n# mutating anything \u2014 text edits below only split existing runs in place, so\n# paragraph indices remain stable across the loop, but each highlight entry\n# must still only be consumed once (two bullets share the same \"Achieved\n# 87% ... 71%\" prefix, differing only by an extra trailing clause).\n$used = @()\nfor ($k = 0; $k -lt $highlights.Count; $k++) {\n    $used += $false\n}\n$paraCount = $d.Paragraphs.Count\n\nfor ($i = 1; $i -le $paraCount; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    $paraText = $para.Range.Text.TrimEnd([char]0x0D, [char]0x07)\n\n    for ($h = 0; $h -lt $highlights.Count; $h++) {\n        if ($used[$h]) { continue }\n        if ($paraText -ne $highlights[$h].Text) { continue }\n        $used[$h] = $true\n\n        foreach ($token in $highlights[$h].Tokens) {\n            $rng = $para.Range\n            $found = $rng.Find.Execute($token)\n            if ($found) {\n                $rng.Font.Bold = $true\n                $rng.Font.Color = $highlightColor\n            }\n        }\n        break\n    }\n}\n"}
